$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 22: replace article #20 content (TOPSMODE / rasagiline paper) ---
# I22 bibliography must come FIRST so the rich-text shared string lands at
# index 15, matching the write order that reproduces the target index layout.
$ws.Range("I22").Value = "Luan, F.; Cordeiro, M.N.; Alonso, N.; García-Mera, X.; Caamaño, O.; Romero-Duran, F.J.; Yañez, M.;  González-Díaz, H. TOPSMODE .model of multiplexing neuroprotective effects of drugs and experimental-theoretic study of new 1,3-rasagiline derivatives potentially useful in neurodegenerative diseases. Bioorg. Med. Chem., 2013, 21(7), 1870-1879. [http://dx.doi.org/10.1016/j.bmc.2013.01. 035] [PMID: 23415089]"
$ws.Range("I22").Font.Size = 10
$ws.Range("I22").Font.Name = "Arial"
$ws.Range("I22").Characters(1,300).Font.Size = 10
$ws.Range("I22").Characters(1,300).Font.Name = "Arial"
$ws.Range("I22").Characters(301,17).Font.Size = 10
$ws.Range("I22").Characters(301,17).Font.Name = "Arial"
$ws.Range("I22").Characters(301,17).Font.Italic = $true
$ws.Range("I22").Characters(318,90).Font.Size = 10
$ws.Range("I22").Characters(318,90).Font.Name = "Arial"
$ws.Range("I22").Font.Italic = $true

$ws.Range("C22").Value = "TOPSMODE`nmodel of multiplexing neuroprotective effects of drugs and experimental- Theoretic study of new 1,3-rasagiline derivatives potentially useful in neurodegenerative diseases."
$ws.Range("D22").Value = 2013
$ws.Range("E22").Value = "Neste estudo foi usado o metodo multiplexing QSAR (mx-QSAR) para ensaios multiplos de molecuals bioativas da basede dados CHEMBL, eles apresentam uma Exatidão de 90%, sensitividade de 98% e seletividade de 80%  para o metodo LDA usado "
$ws.Range("F22").Value = "LDA"
$ws.Range("G22").Value = "CHEMBL"
$ws.Range("H22").Value = "MODESLAB-metodo TOPS-MODE"

# --- Row 21: new article #19 (donepezil-indolyl hybrids paper) ---
$ws.Range("I21").Value = "Bautista-Aguilera, O.M.; Esteban, G.; Bolea, I.; Nikolic, K.; Agbaba,`nD.; Moraleda, I.; Iriepa, I.; Samadi, A.; Soriano, E.; Unzeta,`nM.; Marco-Contelles, J. Design, synthesis, pharmacological evaluation,`nQSAR analysis, molecular modeling and ADMET of novel`ndonepezil-indolyl hybrids as multipotent cholinesterase/monoamine`noxidase inhibitors for the potential treatment of Alzheimer’s disease.`nEur. J. Med. Chem., 2014, 75, 82-95. [http://dx.doi.org/10.`n1016/j.ejmech.2013.12.028] [PMID: 24530494]"
$ws.Range("C21").Value = "Design, synthesis, pharmacological evaluation, QSAR analysis, molecular modeling and ADMET of novel donepezil-indolyl hybrids as multipotent cholinesterase/monoamine oxidase inhibitors for the potential treatment of Alzheimer’s disease"
$ws.Range("D21").Value = 2014
$ws.Range("F21").Value = "PLS"
$ws.Range("E21").Value = " Foram aplicados estudos famacoforicos e 3D-QSAR para desenhar uma serie de novos derivados da donezepila que podem inhibir tanto a AChE como a BChE, aí mostrou-se que tem inhibição no CAS e no PAS das enzimas usando o metodo do PLS.                                                                                                                                                                                  A seleção das variaveis foi feita usando PCA para o maping das estruturas e PLS para fazer a regreção junto com as variaveis descritoras para a valiar os modelos usou-se leave-one-out cross-validation Q2), correlation coefficient (R2 Observed vs. Predicted), Root Main Squared Error of Estimation (RMSEE), and external validation (Root Main Squared Error of Prediction (RMSEP)) "
$ws.Range("G21").Value = "sintesis"

# --- Row heights ---
$ws.Rows(21).RowHeight = 112.2
$ws.Rows(22).RowHeight = 90.6

# --- Column I width ---
$ws.Columns(9).ColumnWidth = 66.21875

# --- Sheet view ---
$ws.Application.ActiveWindow.Zoom = 55
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 5
$ws.Range("E20").Select()

Write-Output "done"
